$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K24").Value = -1.07935118810383
$ws.Range("J25").Value = -1.294239898487671
$ws.Range("K25").Value = -0.004314562288779245
$ws.Range("I26").Value = -1.365431751481502
$ws.Range("J26").Value = -0.07550641528261082
$ws.Range("K26").Value = 2.133689136769128
$ws.Range("H27").Value = -1.391323507574125
$ws.Range("I27").Value = -0.1013981713752329
$ws.Range("J27").Value = 2.107797380676506
$ws.Range("K27").Value = -0.7979924540782646
$ws.Range("G28").Value = -1.557581346555692
$ws.Range("H28").Value = -0.2676560103567999
$ws.Range("I28").Value = 1.941539541694938
$ws.Range("J28").Value = -0.9642502930598316
$ws.Range("K28").Value = -0.6045730858673983
$ws.Range("F29").Value = -1.116218208385255
$ws.Range("G29").Value = 0.1737071278136371
$ws.Range("H29").Value = 2.382902679865375
$ws.Range("I29").Value = -0.5228871548893945
$ws.Range("J29").Value = -0.1632099476969613
$ws.Range("K29").Value = -0.5610502672866085
$ws.Range("E30").Value = -1.220419295757406
$ws.Range("F30").Value = 0.06950604044148601
$ws.Range("G30").Value = 2.278701592493224
$ws.Range("H30").Value = -0.6270882422615456
$ws.Range("I30").Value = -0.2674110350691123
$ws.Range("J30").Value = -0.6652513546587595
$ws.Range("K30").Value = 2.79540342832748
$ws.Range("D31").Value = -1.200984291238337
$ws.Range("E31").Value = 0.08894104496055488
$ws.Range("F31").Value = 2.298136597012293
$ws.Range("G31").Value = -0.6076532377424767
$ws.Range("H31").Value = -0.2479760305500435
$ws.Range("I31").Value = -0.6458163501396906
$ws.Range("J31").Value = 2.814838432846548
$ws.Range("K31").Value = -0.1153445710032998
$ws.Range("C32").Value = -1.28982180544895
$ws.Range("D32").Value = 0.0001035307499418492
$ws.Range("E32").Value = 2.20929908280168
$ws.Range("F32").Value = -0.6964907519530897
$ws.Range("G32").Value = -0.3368135447606565
$ws.Range("H32").Value = -0.7346538643503037
$ws.Range("I32").Value = 2.726000918635935
$ws.Range("J32").Value = -0.2041820852139128
$ws.Range("K32").Value = -0.4407524356782316
$ws.Range("B33").Value = -1.60769821002197
$ws.Range("C33").Value = -0.3177728738230782
$ws.Range("D33").Value = 1.89142267822866
$ws.Range("E33").Value = -1.01436715652611
$ws.Range("F33").Value = -0.6546899493336766
$ws.Range("G33").Value = -1.052530268923324
$ws.Range("H33").Value = 2.408124514062915
$ws.Range("I33").Value = -0.5220584897869329
$ws.Range("J33").Value = -0.7586288402512518
$ws.Range("K33").Value = -0.4330344082701784
$ws.Range("B34").Value = -0.1472629870857898
$ws.Range("C34").Value = 2.061932564965948
$ws.Range("D34").Value = -0.8438572697888214
$ws.Range("E34").Value = -0.4841800625963882
$ws.Range("F34").Value = -0.8820203821860353
$ws.Range("G34").Value = 2.578634400800203
$ws.Range("H34").Value = -0.3515486030496445
$ws.Range("I34").Value = -0.5881189535139633
$ws.Range("J34").Value = -0.26252452153289
$ws.Range("K34").Value = 0.3012009511465751
$ws.Range("B35").Value = 1.996107691246
$ws.Range("C35").Value = -0.9096821435087703
$ws.Range("D35").Value = -0.5500049363163371
$ws.Range("E35").Value = -0.9478452559059842
$ws.Range("F35").Value = 2.512809527080255
$ws.Range("G35").Value = -0.4173734767695933
$ws.Range("H35").Value = -0.6539438272339122
$ws.Range("I35").Value = -0.3283493952528388
$ws.Range("J35").Value = 0.2353760774266263
$ws.Range("K35").Value = -0.05077951757743471
$ws.Range("B36").Value = -1.206907598288802
$ws.Range("C36").Value = -0.8472303910963688
$ws.Range("D36").Value = -1.245070710686016
$ws.Range("E36").Value = 2.215584072300223
$ws.Range("F36").Value = -0.7145989315496251
$ws.Range("G36").Value = -0.9511692820139439
$ws.Range("H36").Value = -0.6255748500328706
$ws.Range("I36").Value = -0.06184937735340548
$ws.Range("J36").Value = -0.3480049723574665
$ws.Range("K36").Value = -0.2261594068364672
$ws.Range("B37").Value = -0.9426906110800261
$ws.Range("C37").Value = -1.340530930669673
$ws.Range("D37").Value = 2.120123852316566
$ws.Range("E37").Value = -0.8100591515332823
$ws.Range("F37").Value = -1.046629501997601
$ws.Range("G37").Value = -0.7210350700165279
$ws.Range("H37").Value = -0.1573095973370627
$ws.Range("I37").Value = -0.4434651923411237
$ws.Range("J37").Value = -0.3216196268201244
$ws.Range("K37").Value = -0.7552991710584124
$ws.Range("B38").Value = -1.239109404396835
$ws.Range("C38").Value = 2.221545378589404
$ws.Range("D38").Value = -0.7086376252604438
$ws.Range("E38").Value = -0.9452079757247627
$ws.Range("F38").Value = -0.6196135437436894
$ws.Range("G38").Value = -0.05588807106422421
$ws.Range("H38").Value = -0.3420436660682852
$ws.Range("I38").Value = -0.2201981005472859
$ws.Range("J38").Value = -0.6538776447855739
$ws.Range("K38").Value = 1.231998594101533
$ws.Range("B39").Value = 2.289204730659292
$ws.Range("C39").Value = -0.6409782731905561
$ws.Range("D39").Value = -0.877548623654875
$ws.Range("E39").Value = -0.5519541916738017
$ws.Range("F39").Value = 0.01177128100566349
$ws.Range("G39").Value = -0.2743843139983975
$ws.Range("H39").Value = -0.1525387484773982
$ws.Range("I39").Value = -0.5862182927156862
$ws.Range("J39").Value = 1.299657946171421
$ws.Range("K39").Value = -0.4039011193250285
$ws.Range("B40").Value = -0.9522025808879372
$ws.Range("C40").Value = -1.188772931352256
$ws.Range("D40").Value = -0.8631784993711827
$ws.Range("E40").Value = -0.2994530266917176
$ws.Range("F40").Value = -0.5856086216957785
$ws.Range("G40").Value = -0.4637630561747793
$ws.Range("H40").Value = -0.8974426004130672
$ws.Range("I40").Value = 0.9884336384740399
$ws.Range("J40").Value = -0.7151254270224096
$ws.Range("K40").Value = -0.7707069585466352
$ws.Range("B41").Value = -1.018399109085777
$ws.Range("C41").Value = -0.6928046771047036
$ws.Range("D41").Value = -0.1290792044252385
$ws.Range("E41").Value = -0.4152347994292995
$ws.Range("F41").Value = -0.2933892339083002
$ws.Range("G41").Value = -0.7270687781465881
$ws.Range("H41").Value = 1.158807460740519
$ws.Range("I41").Value = -0.5447516047559304
$ws.Range("J41").Value = -0.6003331362801561
$ws.Range("K41").Value = 0.1866167262939342
$ws.Range("B42").Value = -0.3420236617375892
$ws.Range("C42").Value = 0.2217018109418759
$ws.Range("D42").Value = -0.0644537840621851
$ws.Range("E42").Value = 0.05739178145881418
$ws.Range("F42").Value = -0.3762877627794737
$ws.Range("G42").Value = 1.509588476107633
$ws.Range("H42").Value = -0.1939705893888161
$ws.Range("I42").Value = -0.2495521209130417
$ws.Range("J42").Value = 0.5373977416610486
$ws.Range("K42").Value = -1.652299669136516
$ws.Range("B43").Value = 0.2693094504600637
$ws.Range("C43").Value = -0.0168461445439973
$ws.Range("D43").Value = 0.104999420977002
$ws.Range("E43").Value = -0.3286801232612859
$ws.Range("F43").Value = 1.557196115625821
$ws.Range("G43").Value = -0.1463629498706283
$ws.Range("H43").Value = -0.2019444813948539
$ws.Range("I43").Value = 0.5850053811792364
$ws.Range("J43").Value = -1.604692029618328
$ws.Range("K43").Value = -0.7667009765463819
$ws.Range("B44").Value = -0.08931394985823571
$ws.Range("C44").Value = 0.03253161566276358
$ws.Range("D44").Value = -0.4011479285755243
$ws.Range("E44").Value = 1.484728310311583
$ws.Range("F44").Value = -0.2188307551848667
$ws.Range("G44").Value = -0.2744122867090923
$ws.Range("H44").Value = 0.5125375758649979
$ws.Range("I44").Value = -1.677159834932566
$ws.Range("J44").Value = -0.8391687818606204
$ws.Range("B45").Value = -0.1051082466437293
$ws.Range("C45").Value = -0.5387877908820172
$ws.Range("D45").Value = 1.34708844800509
$ws.Range("E45").Value = -0.3564706174913596
$ws.Range("F45").Value = -0.4120521490155852
$ws.Range("G45").Value = 0.3748977135585051
$ws.Range("H45").Value = -1.814799697239059
$ws.Range("I45").Value = -0.9768086441671132
$ws.Range("B46").Value = -0.2964732353216546
$ws.Range("C46").Value = 1.589403003565452
$ws.Range("D46").Value = -0.1141560619309969
$ws.Range("E46").Value = -0.1697375934552225
$ws.Range("F46").Value = 0.6172122691188677
$ws.Range("G46").Value = -1.572485141678696
$ws.Range("H46").Value = -0.7344940886067506
$ws.Range("B47").Value = 1.553990447017354
$ws.Range("C47").Value = -0.1495686184790955
$ws.Range("D47").Value = -0.2051501500033212
$ws.Range("E47").Value = 0.5817997125707691
$ws.Range("F47").Value = -1.607897698226795
$ws.Range("G47").Value = -0.7699066451548492
$ws.Range("B48").Value = -0.3583240076259202
$ws.Range("C48").Value = -0.4139055391501458
$ws.Range("D48").Value = 0.3730443234239444
$ws.Range("E48").Value = -1.81665308737362
$ws.Range("F48").Value = -0.9786620343016739
$ws.Range("B49").Value = -0.2364856966075551
$ws.Range("C49").Value = 0.5504641659665352
$ws.Range("D49").Value = -1.639233244831029
$ws.Range("E49").Value = -0.8012421917590831
$ws.Range("B50").Value = 0.6235572686657249
$ws.Range("C50").Value = -1.566140142131839
$ws.Range("D50").Value = -0.7281490890598934
$ws.Range("B51").Value = -1.668261113776646
$ws.Range("C51").Value = -0.8302700607047004
$ws.Range("B52").Value = -0.7276043929666616
